# Add the "Department Code" lookup sheet after Sheet1, populate it, and
# move the Sheet1 selection cursor, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Department Code"

$ws2.Cells.Item(1, 1).Value = "Department"
$ws2.Cells.Item(1, 2).Value = "Code"
$ws2.Cells.Item(2, 1).Value = "HR - BCD"
$ws2.Cells.Item(2, 2).Value = "HRB"
$ws2.Cells.Item(3, 1).Value = "Billing Department"
$ws2.Cells.Item(3, 2).Value = "BIL"
$ws2.Cells.Item(4, 1).Value = "EIC"
$ws2.Cells.Item(4, 2).Value = "EIC"
$ws2.Cells.Item(5, 1).Value = "Environment/PCO"
$ws2.Cells.Item(5, 2).Value = "ENV"
$ws2.Cells.Item(6, 1).Value = "Finance Department"
$ws2.Cells.Item(6, 2).Value = "FIN"
$ws2.Cells.Item(7, 1).Value = "Fuel and Lube Management"
$ws2.Cells.Item(7, 2).Value = "FLM"
$ws2.Cells.Item(8, 1).Value = "Health and Safety"
$ws2.Cells.Item(8, 2).Value = "HAS"
$ws2.Cells.Item(9, 1).Value = "IT Department - BCD"
$ws2.Cells.Item(9, 2).Value = "ITB"
$ws2.Cells.Item(10, 1).Value = "IT Department - SITE"
$ws2.Cells.Item(10, 2).Value = "ITS"
$ws2.Cells.Item(11, 1).Value = "Laboratory and Chemical"
$ws2.Cells.Item(11, 2).Value = "LAB"
$ws2.Cells.Item(12, 1).Value = "Maintenance"
$ws2.Cells.Item(12, 2).Value = "MAI"
$ws2.Cells.Item(13, 1).Value = "Office of the GM"
$ws2.Cells.Item(13, 2).Value = "OOG"
$ws2.Cells.Item(14, 1).Value = "Operation"
$ws2.Cells.Item(14, 2).Value = "OPE"
$ws2.Cells.Item(15, 1).Value = "Purchasing Department"
$ws2.Cells.Item(15, 2).Value = "PUR"
$ws2.Cells.Item(16, 1).Value = "Reconditioning"
$ws2.Cells.Item(16, 2).Value = "REC"
$ws2.Cells.Item(17, 1).Value = "Security"
$ws2.Cells.Item(17, 2).Value = "SEC"
$ws2.Cells.Item(18, 1).Value = "HR - SITE"
$ws2.Cells.Item(18, 2).Value = "HRS"
$ws2.Cells.Item(19, 1).Value = "Special Proj/Facilities Imp"
$ws2.Cells.Item(19, 2).Value = "SPE"
$ws2.Cells.Item(20, 1).Value = "Trading Department"
$ws2.Cells.Item(20, 2).Value = "TRA"
$ws2.Cells.Item(21, 1).Value = "Warehouse - CENPRI "
$ws2.Cells.Item(21, 2).Value = "WHC"
$ws2.Cells.Item(22, 1).Value = "Warehouse - Progen "
$ws2.Cells.Item(22, 2).Value = "WHP"

# Header row is bold (reuses the workbook's existing bold style).
$ws2.Range("A1:B1").Font.Bold = $true

# Column A sized to fit the longest department name.
$ws2.Columns.Item(1).ColumnWidth = 25.166666666666668

# Restore selections: Department Code sheet was left on E8, Sheet1 stays
# the active tab with N7 selected.
$ws2.Range("E8").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("N7").Select() | Out-Null

